$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F with schedule availability notes, mirroring the layout of column B
$ws.Range("F3").Value = "a partir das 19:30"
$ws.Range("F4").Value = "a partir das 19:30"
$ws.Range("F5").Value = "a partir das 19:30"
$ws.Range("F6").Value = "a partir das 19:30"
$ws.Range("F7").Value = "a partir das 18:30"
$ws.Range("F8").Value = "até o meio dia"
$ws.Range("F9").Value = "**"

# Widen column F so the new text fits
$ws.Columns.Item(6).ColumnWidth = 18.67

# Move/update the active selection to the last cell that was edited
$ws.Range("F9").Select() | Out-Null
